# Fix the misspelled product name "Blackly"/"Blocky" -> "Blockly" on the
# "code generation feature" slide (slide 5), keeping existing run
# formatting (lang/altLang/dirty/kumimoji/smtClean) intact.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# First paragraph: 「Blackly」 -> 「Blockly」
$hit = $tr.Find("Blackly", 0, $false, $false)
if ($hit) {
    $hit.Text = "Blockly"
}

# Second paragraph: 「Blocky」 -> 「Blockly」
$hit2 = $tr.Find("Blocky", 0, $false, $false)
if ($hit2) {
    $hit2.Text = "Blockly"
}
